$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($cellAddr, $text) {
    $c = $ws.Range($cellAddr)
    $c.Value = "'" + $text
    $c.Style = "Normal"
}

Set-TextValue "D2" "313.43"
Set-TextValue "E2" "5.70%"
Set-TextValue "D3" "44.69"
Set-TextValue "E3" "7.03%"
Set-TextValue "D4" "5.131"
Set-TextValue "E4" "1.81%"
Set-TextValue "D5" "0.08044"
Set-TextValue "E5" "6.49%"
Set-TextValue "E6" "2.90%"
Set-TextValue "D7" "1.693"
Set-TextValue "E7" "6.15%"
Set-TextValue "D8" "1.083"
Set-TextValue "E8" "16.63%"
Set-TextValue "D9" "0.1296"
Set-TextValue "E9" "8.50%"
Set-TextValue "D10" "0.1918"
Set-TextValue "E10" "4.30%"
Set-TextValue "D11" "0.09392"
Set-TextValue "E11" "4.45%"
Set-TextValue "D12" "0.04233"
Set-TextValue "E12" "6.28%"
Set-TextValue "D13" "0.1042"
Set-TextValue "E13" "-0.93%"
Set-TextValue "D14" "0.001311"
Set-TextValue "E14" "2.26%"
Set-TextValue "D15" "0.005922"
Set-TextValue "E15" "1.58%"
Set-TextValue "D17" "3.392"
Set-TextValue "E17" "1.14%"
Set-TextValue "D18" "2.403"
Set-TextValue "E18" "-0.23%"
Set-TextValue "D19" "0.3373"
Set-TextValue "E19" "1.60%"
Set-TextValue "D20" "8.051"
Set-TextValue "E20" "2.06%"
Set-TextValue "E21" "-3.50%"
Set-TextValue "D23" "0.04204"
Set-TextValue "E23" "3.63%"
Set-TextValue "D24" "0.001271"
Set-TextValue "E24" "0.38%"
Set-TextValue "D25" "0.004566"
Set-TextValue "E25" "15.92%"
Set-TextValue "D26" "0.0001341"
Set-TextValue "E26" "8.97%"
Set-TextValue "D38" "0.02715"
Set-TextValue "E38" "12.52%"
Set-TextValue "D39" "0.05423"
Set-TextValue "E39" "3.95%"
Set-TextValue "D40" "0.005635"
Set-TextValue "E40" "-6.63%"
Set-TextValue "D41" "0.007743"
Set-TextValue "E41" "-0.41%"
Set-TextValue "D42" "0.1414"
Set-TextValue "E42" "6.32%"
Set-TextValue "D43" "0.007336"
Set-TextValue "E43" "-2.67%"
Set-TextValue "D44" "0.007948"
Set-TextValue "E44" "1.23%"
Set-TextValue "D45" "0.3130"
Set-TextValue "E45" "-2.82%"
Set-TextValue "D46" "0.00006793"
Set-TextValue "E46" "0.11%"
Set-TextValue "E47" "-0.70%"
Set-TextValue "D48" "0.06955"
Set-TextValue "E48" "54.64%"
Set-TextValue "D49" "0.003975"
Set-TextValue "E49" "-5.42%"
Set-TextValue "D50" "0.00002087"
Set-TextValue "E50" "-0.70%"
Set-TextValue "E51" "-0.70%"
